$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45922
$ws.Range("B2").Value = 30.9
$ws.Range("C2").Value = 30.9
$ws.Range("D2").Value = 26.2
$ws.Range("E2").Value = 26.2
$ws.Range("F2").Value = 34.9
$ws.Range("G2").Value = 30.14
$ws.Range("H2").Value = 30.9
$ws.Range("I2").Value = 60.6
$ws.Range("J2").Value = 62.03
$ws.Range("K2").Value = 29.35
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0.01
$ws.Range("P2").Value = -0.01
$ws.Range("Q2").Value = -0.01
$ws.Range("R2").Value = -0.01
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 15.6
$ws.Range("U2").Value = 67.98999999999999
$ws.Range("V2").Value = 108.95
$ws.Range("W2").Value = 109.38
$ws.Range("X2").Value = 82.8
$ws.Range("Y2").Value = 80
$ws.Range("Z2").Value = 34.53
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 95.28
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 109.16
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 81.40000000000001
$ws.Range("AG2").Value = "0h-18h"
